$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns E1:I1
$ws.Range("E1").Value = "i5"
$ws.Range("F1").Value = "i6"
$ws.Range("G1").Value = "i7"
$ws.Range("H1").Value = "i8"
$ws.Range("I1").Value = "i9"

# Data row 2 - update existing values and add new columns
$ws.Range("A2").Value = 0.001183285792159792
$ws.Range("B2").Value = 0.4613257385728005
$ws.Range("C2").Value = 0.2240924011402907
$ws.Range("D2").Value = 0.5809271786130042
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
